# Fruta / hortaliza, semanal
#
# Adds the latest weekly price observation for "Coco" (Mercado Mayorista Lo
# Valledor de Santiago) by inserting a new row right above the current
# row 36, shifting all the following rows down by one (dimension grows
# from A1:T66 to A1:T67), and populating the new row with the new week's
# data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36 - pushes existing rows 36..66 down to 37..67
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record
$ws.Cells.Item(36, 1).Value  = 6
$ws.Cells.Item(36, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(36, 3).Value  = "Metropolitana"
$ws.Cells.Item(36, 4).Value  = 44719
$ws.Cells.Item(36, 5).Value  = 13
$ws.Cells.Item(36, 6).Value  = "Fruta"
$ws.Cells.Item(36, 7).Value  = 100108
$ws.Cells.Item(36, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(36, 9).Value  = 100108007
$ws.Cells.Item(36, 10).Value = "Coco"
$ws.Cells.Item(36, 11).Value = "Sin especificar"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 200
$ws.Cells.Item(36, 14).Value = 22000
$ws.Cells.Item(36, 15).Value = 23000
$ws.Cells.Item(36, 16).Value = 22500
$ws.Cells.Item(36, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(36, 18).Value = "Perú"
$ws.Cells.Item(36, 19).Value = 1125
$ws.Cells.Item(36, 20).Value = 20
